$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.026.76'
$ws.Range("E2").Value = '  -1.83%  '
$ws.Range("D3").Value = '2.103.04'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("E4").Value = '  -0.68%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '349.54'
$ws.Range("E5").Value = '  +3.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5154'
$ws.Range("E7").Value = '  -1.76%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4430'
$ws.Range("E8").Value = '  -2.78%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.38'
$ws.Range("E9").Value = '  -4.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08968'
$ws.Range("E10").Value = '  -1.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.168'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.44'
$ws.Range("E12").Value = '  +3.88%  '
$ws.Range("D13").Value = '2.101.14'
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.213'
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.729'
$ws.Range("E15").Value = '  -1.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '98.91'
$ws.Range("E16").Value = '  +1.89%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001146'
$ws.Range("E17").Value = '  -2.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.005'
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.81'
$ws.Range("E19").Value = '  +7.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06672'
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("E21").Value = '  -0.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.209'
$ws.Range("E22").Value = '  -1.56%  '
$ws.Range("D23").Value = '30.123.23'
$ws.Range("E23").Value = '  -1.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.65'
$ws.Range("E24").Value = '  -2.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.339'
$ws.Range("E25").Value = '  -0.56%  '
$ws.Range("D26").Value = '2.351.74'
$ws.Range("E26").Value = '  -0.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.95'
$ws.Range("E27").Value = '  -1.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.550'
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '161.89'
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.30'
$ws.Range("E30").Value = '  -0.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.172'
$ws.Range("E31").Value = '  -3.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1060'
$ws.Range("E32").Value = '  -1.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.641'
$ws.Range("E33").Value = '  -0.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.211'
$ws.Range("E34").Value = '  -2.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.970'
$ws.Range("E35").Value = '  +0.73%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.959'
$ws.Range("E36").Value = '  +1.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.16'
$ws.Range("E37").Value = '  -4.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02569'
$ws.Range("E38").Value = '  -2.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06787'
$ws.Range("E39").Value = '  -0.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2284'
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.345'
$ws.Range("E41").Value = '  +7.17%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.49'
$ws.Range("E42").Value = '  -0.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6786'
$ws.Range("E43").Value = '  -1.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.21'
$ws.Range("E44").Value = '  -3.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6358'
$ws.Range("E45").Value = '  -1.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.282'
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000361'
$ws.Range("E47").Value = '  -2.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.643'
$ws.Range("E48").Value = '  -1.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.217'
$ws.Range("E49").Value = '  -3.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '82.21'
$ws.Range("E50").Value = '  -1.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07225'
$ws.Range("E51").Value = '  +0.10%  '
